$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update the "through" date references from Oct 25 -> Oct 26
$ws.Name = "Through 2022-10-26"
$ws.Range("B1").Value = "October 2022 (through October 26)"

# --- Numeric cell updates (new daily carjacking data for 2022-11-03) ---

# Row 2 - Garfield Park
$ws.Range("B2").Value = 5
$ws.Range("L2").Value = 18
$ws.Range("V2").Value = 17
$ws.Range("AP2").Value = 7

# Row 5 - North Lawndale
$ws.Range("L5").Value = 15
$ws.Range("AF5").Value = 1
$ws.Range("BJ5").Value = 2

# Row 6 - Austin
$ws.Range("B6").Value = 2
$ws.Range("AZ6").Value = 7
$ws.Range("BJ6").Value = 3
$ws.Range("BT6").Value = 3

# Row 7 - Englewood
$ws.Range("L7").Value = 5

# Row 10 - New City
$ws.Range("B10").Value = 5
$ws.Range("AF10").Value = 2

# Row 18 - Grand Crossing
$ws.Range("BT18").Value = 1

# Row 26 - Lake View
$ws.Range("AF26").Value = 2

# Row 30 - West Town
$ws.Range("L30").Value = 8

# Row 98 - Woodlawn
$ws.Range("B98").Value = 1
$ws.Range("BJ98").Value = 2
